# Update latest output (run 147)
$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
$wsSchedule.Range("E3").Value = -191.23649025
$wsSchedule.Range("F3").Value = -5.621296009700176

$wsSchedule.Range("E4").Value = 631.528716
$wsSchedule.Range("F4").Value = 30.37656161616162

$wsSchedule.Range("E5").Value = -14.7904185
$wsSchedule.Range("F5").Value = -0.4118746449456977

# --- Detailed sheet updates ---
$wsDetailed.Range("B30").Value = -27
$wsDetailed.Range("B31").Value = -25.94511

$wsDetailed.Range("B32").Value = -9.85553
$wsDetailed.Range("C32").Value = "historical"

$wsDetailed.Range("B33").Value = -5.71179
$wsDetailed.Range("C33").Value = "historical"

$wsDetailed.Range("B34").Value = -5.01
$wsDetailed.Range("B35").Value = 0
$wsDetailed.Range("B36").Value = 36.06018
$wsDetailed.Range("B37").Value = 48.95598
$wsDetailed.Range("B38").Value = 57.31
$wsDetailed.Range("B39").Value = 66.09419
$wsDetailed.Range("B41").Value = 73.2
$wsDetailed.Range("B42").Value = 75.07599
$wsDetailed.Range("B43").Value = 73.19
$wsDetailed.Range("B45").Value = 58.65096
$wsDetailed.Range("B47").Value = 57.31
$wsDetailed.Range("B49").Value = 61.30589
$wsDetailed.Range("B51").Value = 62.20481
$wsDetailed.Range("B53").Value = 57.06
$wsDetailed.Range("B55").Value = 57.06
$wsDetailed.Range("B56").Value = 63.71455
$wsDetailed.Range("B57").Value = 63.33347
$wsDetailed.Range("B58").Value = 63.32307
$wsDetailed.Range("B59").Value = 57.09
$wsDetailed.Range("B60").Value = 57.31
$wsDetailed.Range("B61").Value = 64.89
$wsDetailed.Range("B62").Value = 62.81908
$wsDetailed.Range("B63").Value = 51.93597
$wsDetailed.Range("B64").Value = 35.88
$wsDetailed.Range("B66").Value = -5.1403
$wsDetailed.Range("B67").Value = -7.61355
$wsDetailed.Range("B69").Value = -9.67783
$wsDetailed.Range("B71").Value = -7.71363
$wsDetailed.Range("B72").Value = -7.85575
$wsDetailed.Range("B74").Value = -7.82081
$wsDetailed.Range("B75").Value = -7.89798
$wsDetailed.Range("B76").Value = -7.6582
$wsDetailed.Range("B77").Value = -5.71213
$wsDetailed.Range("B78").Value = -5.55375
$wsDetailed.Range("B79").Value = -5.17224
$wsDetailed.Range("B80").Value = -5.01
$wsDetailed.Range("B82").Value = -1.17711
$wsDetailed.Range("B83").Value = -2.95369
$wsDetailed.Range("B84").Value = -11.01
$wsDetailed.Range("B85").Value = -8.41559
$wsDetailed.Range("B87").Value = 0
$wsDetailed.Range("B88").Value = 5.32284
$wsDetailed.Range("B89").Value = 59.36649
$wsDetailed.Range("B90").Value = 53.90468
$wsDetailed.Range("B91").Value = 43.36402
$wsDetailed.Range("B92").Value = 38.26315
$wsDetailed.Range("B93").Value = 43.57993
$wsDetailed.Range("B94").Value = 30.46316
$wsDetailed.Range("B95").Value = 56.52657
$wsDetailed.Range("B96").Value = 56.85286
$wsDetailed.Range("B97").Value = 52.45737
